# Apply the "Generate Yearly Report" configuration update to the REFramework
# Data/Config.xlsx workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Row 2 - Orchestrator queue name (no longer bold-styled)
$settings.Range("B2").Value = "WI4_Items"
$settings.Range("B2").Style = "Normal"

# Row 3 - Orchestrator queue folder
$settings.Range("B3").Value = "Shared"

# Row 5 - Business process name (logging)
$settings.Range("B5").Value = "Generate Yearly Report"

# Row 7 - System1 URL (new setting)
$settings.Range("A7").Value = "System1Url"
$settings.Range("B7").Value = "https://acme-test.uipath.com"

# Row 9 - System1 credential (new setting)
$settings.Range("A9").Value = "System1Credential"
$settings.Range("B9").Value = "System1Credential"

# Row 11 - Work items URL (new setting)
$settings.Range("A11").Value = "WorkItemsUrl"
$settings.Range("B11").Value = "https://acme-test.uipath.com/work-items"

# Row 13 - Process name (new setting)
$settings.Range("A13").Value = "ProcessName"
$settings.Range("B13").Value = "chrome"

# Row 15 - Reset data URL (new setting)
$settings.Range("A15").Value = "ResetDataUrl"
$settings.Range("B15").Value = "https://acme-test.uipath.com/reset-test-data"

# Row 17 - Download monthly reports URL (new setting)
$settings.Range("A17").Value = "DownloadMonthlyReportsUrl"
$settings.Range("B17").Value = "https://acme-test.uipath.com/reports/download"

# Row 19 - Upload yearly reports URL (new setting)
$settings.Range("A19").Value = "UploadYearlyReportsUrl"
$settings.Range("B19").Value = "https://acme-test.uipath.com/reports/upload"

# Row 21 - Status (new setting)
$settings.Range("A21").Value = "Status"
$settings.Range("B21").Value = "Completed"

# Row 23 - Year (new setting)
$settings.Range("A23").Value = "Year"
$settings.Range("B23").Value = 2022
$settings.Range("B23").HorizontalAlignment = -4131

# Row 25 - Report directory (new setting)
$settings.Range("A25").Value = "Report_Dir"
$settings.Range("B25").Value = "Data\Report\"

# Row 27 - Timeout download (new setting)
$settings.Range("A27").Value = "TimeOutDownload"
$settings.Range("B27").Value = 10
$settings.Range("B27").HorizontalAlignment = -4131

# Row heights recalculated for wrapped description cells
$settings.Rows.Item(3).RowHeight = 45
$settings.Rows.Item(5).RowHeight = 30

# Remove the three trailing blank rows (996-998) no longer present in the sheet
$settings.Rows.Item(996).Delete()
$settings.Rows.Item(996).Delete()
$settings.Rows.Item(996).Delete()

# Active selection on the Settings sheet
$settings.Range("C23").Select()

# ---------------------------------------------------------------------
# Constants sheet - row heights match the same re-flow as Settings
# ---------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")
$constants.Rows.Item(2).RowHeight = 30
$constants.Rows.Item(3).RowHeight = 45
$constants.Rows.Item(17).RowHeight = 45
